# Updates cryptocurrency price/volume figures (and restores the correct
# MantraDAO/ImmutableX row order) on Sheet1, per the Nov 26 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "92.225.19"
$ws.Range("E2").Value = "  -1.27%  "

# Row 3
$ws.Range("D3").Value = "3.338.62"
$ws.Range("E3").Value = "  -1.77%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.07%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "615.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.85%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.41"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.55%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.386"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.15%  "

# Row 9
$ws.Range("E9").Value = "  -0.10%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.956"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.54%  "

# Row 11
$ws.Range("D11").Value = "3.337.71"
$ws.Range("E11").Value = "  -1.66%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.32%  "

# Row 13
$ws.Range("E13").Value = "  -0.75%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.39%  "

# Row 15
$ws.Range("D15").Value = "92.012.10"
$ws.Range("E15").Value = "  -1.49%  "

# Row 16
$ws.Range("D16").Value = "3.966.58"
$ws.Range("E16").Value = "  -1.77%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000243"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.63%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.90%  "

# Row 19
$ws.Range("D19").Value = "3.336.19"
$ws.Range("E19").Value = "  -1.91%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.21%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.55%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.80%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "495.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.80%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.445"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.66%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000183"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.84%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "89.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.64%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.54%  "

# Row 29
$ws.Range("D29").Value = "3.512.58"
$ws.Range("E29").Value = "  -1.88%  "

# Row 30
$ws.Range("E30").Value = "  +0.69%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.139"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.77%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.22%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.990"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.94%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.172"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.80%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.61%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.528"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.73%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "565.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.84%  "

# Row 39
$ws.Range("E39").Value = "  +0.01%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.82%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.148"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.33%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.872"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.68%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.41%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0415"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.08%  "

# Row 46
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "

# Row 47
$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.66%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.25%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.98%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.07%  "

